$wb = $excel.ActiveWorkbook

# "想去人数" (want-to-go count) refresh for two events, mirrored on both the
# "展览" sheet and the "全部类型" aggregate sheet.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 1411
    $ws.Range("F6").Value = 7
}
